$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking strings
# (e.g. "1.001", "312.20") are stored as literal text, not converted to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '28.233.30'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '1.874.19'
$ws.Range('E3').Value = '  +4.13%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = '312.20'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '0.5028'
$ws.Range('E7').Value = '  -1.38%  '
$ws.Range('D8').Value = '0.3938'
$ws.Range('E8').Value = '  +0.89%  '
$ws.Range('D9').Value = '0.09859'
$ws.Range('E9').Value = '  +26.94%  '
$ws.Range('D10').Value = '1.142'
$ws.Range('E10').Value = '  +3.79%  '
$ws.Range('D11').Value = '41.22'
$ws.Range('D12').Value = '6.482'
$ws.Range('E12').Value = '  +2.49%  '
$ws.Range('D13').Value = '21.03'
$ws.Range('E13').Value = '  +4.02%  '
$ws.Range('D14').Value = '1.868.91'
$ws.Range('E14').Value = '  +3.88%  '
$ws.Range('D15').Value = '1.000'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '7.408'
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').Value = '0.00001136'
$ws.Range('E17').Value = '  +6.08%  '
$ws.Range('D18').Value = '93.54'
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').Value = '0.06630'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').Value = '17.45'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').Value = '6.147'
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('D23').Value = '28.290.89'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').Value = '11.35'
$ws.Range('E24').Value = '  +2.59%  '
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('E26').Value = '  +5.79%  '
$ws.Range('D27').Value = '21.33'
$ws.Range('E27').Value = '  +4.83%  '
$ws.Range('D28').Value = '2.085.09'
$ws.Range('E28').Value = '  +3.65%  '
$ws.Range('D29').Value = '158.07'
$ws.Range('D30').Value = '127.59'
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('D31').Value = '0.1067'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('D32').Value = '1.066'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').Value = '5.635'
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('D34').Value = '3.620'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('D35').Value = '0.06812'
$ws.Range('E35').Value = '  -3.07%  '
$ws.Range('D36').Value = '9.542'
$ws.Range('E36').Value = '  +5.05%  '
$ws.Range('D37').Value = '0.02393'
$ws.Range('E37').Value = '  +2.39%  '
$ws.Range('D38').Value = '0.2190'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '5.028'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '11.52'
$ws.Range('E40').Value = '  +0.41%  '
$ws.Range('D41').Value = '0.6315'
$ws.Range('E41').Value = '  +3.21%  '
$ws.Range('D42').Value = '1.172'
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('D44').Value = '13.61'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('D45').Value = '0.6031'
$ws.Range('E45').Value = '  +2.27%  '
$ws.Range('D46').Value = '3.668'
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').Value = '1.270'
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').Value = '124.88'
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('D49').Value = '1.995'
$ws.Range('E49').Value = '  +4.53%  '
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').Value = '1.124'
$ws.Range('E51').Value = '  +6.18%  '

# Restore column D formatting/style back to the workbook default so no
# extra per-cell style index is left behind on the data cells.
$dRange.Style = "Normal"
